$wb = $excel.ActiveWorkbook

# --- Add the new "Hoja2" worksheet after DIGITAL_ASSETS ---
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Hoja2"
$digitalAssets = $wb.Worksheets.Item("DIGITAL_ASSETS")
$newSheet.Move($null, $digitalAssets)

# Re-fetch a live reference to the moved sheet
$ws2 = $wb.Worksheets.Item("Hoja2")

# --- Populate Hoja2 with the SKU / technical-document filename matrix ---
$data = @(
    @(246021151, "206191~1.PDF", 23.85),
    @(246021181, "246021151-piso-pared-nuevo-tahoe-azul-mt-ft.pdf", 57.45),
    @(246021451, "246021181-piso-pared-nuevo-tahoe-azul-oscuro-mt-ft.pdf", 68.65),
    @(246021491, "246021451-piso-pared-nuevo-tahoe-verde-mt-ft.pdf", 68.65),
    @(247031151, "246021491-piso-pared-nuevo-tahoe-verde-oscuro-mt-ft.pdf", 79.85),
    @(247031761, "247031151-piso-pared-tulum-azul-mt-ft.pdf", 46.25),
    @(247041101, "247031761-piso-pared-tulum-cafe-mt-ft.pdf", 46.25)
)

$r = 1
foreach ($row in $data) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Range($ws2.Cells.Item($r, 1), $ws2.Cells.Item($r, 2)).WrapText = $true
    $ws2.Rows.Item($r).RowHeight = $row[2]
    $r = $r + 1
}

# Last row: SKU cell left empty, only the filename is populated
$ws2.Cells.Item(8, 2).Value = "247041101-piso-pared-oceano-azul-claro-cd-ft.pdf"
$ws2.Range($ws2.Cells.Item(8, 1), $ws2.Cells.Item(8, 2)).WrapText = $true
$ws2.Rows.Item(8).RowHeight = 57.45

$ws2.Range("A1").Select()

# --- Format the DIGITAL_ASSETS sheet: header row + body get "locked" styling ---
$ws1 = $wb.Worksheets.Item("DIGITAL_ASSETS")
$usedRange = $ws1.Range("A2:B109")
$usedRange.Locked = $true

$header = $ws1.Range("A1:B1")
$header.Locked = $true

# --- Make Hoja2 the active/visible tab ---
$ws1.Range("A10").Select()
$ws2.Activate()
$ws2.Range("A1").Select()
